# Update "contadores" counts per CI parser fix (diff: results/contadores.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Jenkins
$ws.Range("B2").Value = 2

# Row 3 - Travis
$ws.Range("B3").Value = 31
$ws.Range("C3").Value = 0

# Row 4 - Circle CI
$ws.Range("B4").Value = 19
$ws.Range("C4").Value = 0

# Row 5 - GitHub Actions
$ws.Range("B5").Value = 130
$ws.Range("C5").Value = 0

# Row 6 - Azure Pipelines
$ws.Range("B6").Value = 4

# Row 7 - Bamboo
$ws.Range("B7").Value = 0

# Row 9 - GitLab CI
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 0

# Row 12 - Bazel
$ws.Range("B12").Value = 2

# Row 15 - Totales
$ws.Range("B15").Value = 153
$ws.Range("C15").Value = 0
